# New crime data collected - weekly CompStat update for 123rd Precinct
# (covering week 3/27/2023 through 4/2/2023, Volume 30 Number 13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings) -------------------------------
# ---------------------------------------------------------------------
# A8  : "Volume 30   Number  12"  -> "...  13"
# C9  : "Report Covering the Week  3/20/2023  Through  3/26/2023"
#        -> "...3/27/2023  Through  4/2/2023"
$ws.Range("A8").Value = "Volume 30   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/27/2023  Through  4/2/2023"

# ---------------------------------------------------------------------
# Helper: set a cell to the literal text "0" (used for a zero count
# that is displayed as text rather than as a formatted number) while
# re-using the workbook's existing "text" style for that row so the
# look (font/alignment) matches the surrounding text cells exactly.
# ---------------------------------------------------------------------
function Set-TextZero($cellRef, $formatSourceRef) {
    $src = $ws.Range($formatSourceRef)
    $dst = $ws.Range($cellRef)
    $dst.NumberFormat = "@"
    $dst.Value = "0"
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

function Set-TextNA($cellRef, $formatSourceRef) {
    $src = $ws.Range($formatSourceRef)
    $dst = $ws.Range($cellRef)
    $dst.NumberFormat = "@"
    $dst.Value = "***.*"
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-TextZero "F15" "G15"

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("N16").Value = -76.923076923076

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 700
$ws.Range("I17").Value = 20
$ws.Range("K17").Value = 81.818181818181
$ws.Range("L17").Value = 185.714285714286
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = 25

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Set-TextZero "C18" "D18"
$ws.Range("F18").Value = 6
$ws.Range("M18").Value = -26.923076923076
$ws.Range("N18").Value = -70.3125

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 5
Set-TextZero "D19" "C18"
Set-TextNA "E19" "E18"
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 67
$ws.Range("K19").Value = -6.944444444444
$ws.Range("L19").Value = 52.272727272727
$ws.Range("M19").Value = 67.5
$ws.Range("N19").Value = 97.058823529411

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-TextZero "D20" "C20"
Set-TextNA "E20" "E18"
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("N20").Value = -89.655172413793

# ---------------------------------------------------------------------
# Row 21 - TOTAL (bold row)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 7
Set-TextZero "D21" "A21"
Set-TextNA "E21" "A21"
$ws.Range("F21").Value = 41
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 36.666666666666
$ws.Range("I21").Value = 128
$ws.Range("K21").Value = 13.274336283185
$ws.Range("L21").Value = 91.044776119403
$ws.Range("M21").Value = 40.659340659340
$ws.Range("N21").Value = -57.894736842105

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 33
$ws.Range("H24").Value = 33.333333333333
$ws.Range("I24").Value = 136
$ws.Range("J24").Value = 99
$ws.Range("K24").Value = 37.373737373737
$ws.Range("L24").Value = 151.851851851852
$ws.Range("M24").Value = 22.522522522522

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -35
$ws.Range("I25").Value = 50
$ws.Range("J25").Value = 49
$ws.Range("K25").Value = 2.040816326530
$ws.Range("L25").Value = 92.307692307692
$ws.Range("M25").Value = 6.382978723404

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
Set-TextZero "F26" "G26"

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -57.142857142857
